# Insert a new weekly price record as row 3 (right after the header data
# already present in row 2), pushing all existing data rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (old 3..64) down to (4..65) by inserting a new
# row at position 3.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new record.
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 44699
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 100112029
$ws.Range("G3").Value = "Orégano"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 16
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("N3").Value = "$/docena de atados"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 5000
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = "Hortaliza"
